# Updates cryptos list values (price + volume) per the sync diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.247.50"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.969.30"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.59"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.59"
$ws.Range("E7").Value = "  -9.30%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.01"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.857"
$ws.Range("E13").Value = "  -8.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.41"
$ws.Range("E14").Value = "  +8.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.03"
$ws.Range("E15").Value = "  -7.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.255.56"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.968.33"
$ws.Range("E18").Value = "  -3.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.129.23"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.24"
$ws.Range("E20").Value = "  -3.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.20"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -5.91%  "
$ws.Range("E26").Value = "  -4.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.23"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.86"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  +10.18%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -7.29%  "
$ws.Range("E33").Value = "  -5.91%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  -7.87%  "
$ws.Range("E36").Value = "  +4.73%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.09"
$ws.Range("E40").Value = "  +9.45%  "
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.22"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.31"
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.97"
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.52"
$ws.Range("E48").Value = "  -7.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.339.07"
$ws.Range("E49").Value = "  -6.42%  "
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.149.70"
$ws.Range("E51").Value = "  -3.91%  "
